$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "26.007.04"
$cell.Style = $origStyle

$ws.Range("E2").Value = "  +0.23%  "
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.641.10"
$cell.Style = $origStyle

$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.07%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "215.22"
$cell.Style = $origStyle

$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.255"
$cell.Style = $origStyle

$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0637"
$cell.Style = $origStyle

$ws.Range("E9").Value = "  +0.03%  "
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "19.55"
$cell.Style = $origStyle

$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("E12").Value = "  +0.08%  "
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.626.99"
$cell.Style = $origStyle

$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("E14").Value = "  -0.22%  "
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "63.41"
$cell.Style = $origStyle

$ws.Range("E15").Value = "  +1.35%  "
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0759"
$cell.Style = $origStyle

$ws.Range("E16").Value = "  -0.05%  "
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "26.049.47"
$cell.Style = $origStyle

$ws.Range("E18").Value = "  +0.18%  "
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "194.41"
$cell.Style = $origStyle

$ws.Range("E19").Value = "  +0.13%  "
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.35"
$cell.Style = $origStyle

$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("E21").Value = "  -0.40%  "
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.19"
$cell.Style = $origStyle

$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("E23").Value = "  +4.16%  "
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "143.83"
$cell.Style = $origStyle

$ws.Range("E24").Value = "  -0.12%  "
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.78"
$cell.Style = $origStyle

$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  +0.30%  "
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0494"
$cell.Style = $origStyle

$ws.Range("E30").Value = "  -1.32%  "
$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.25"
$cell.Style = $origStyle

$ws.Range("E31").Value = "  +0.70%  "
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.27"
$cell.Style = $origStyle

$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("E33").Value = "  -0.20%  "
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.45"
$cell.Style = $origStyle

$ws.Range("E34").Value = "  +0.88%  "
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.902"
$cell.Style = $origStyle

$ws.Range("E35").Value = "  -0.20%  "
$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.128.17"
$cell.Style = $origStyle

$ws.Range("E36").Value = "  -0.97%  "
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.537"
$cell.Style = $origStyle

$ws.Range("E37").Value = "  -1.47%  "
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.45"
$cell.Style = $origStyle

$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E39").Value = "  -0.23%  "
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "98.69"
$cell.Style = $origStyle

$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("E41").Value = "  +0.18%  "
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.795"
$cell.Style = $origStyle

$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("E43").Value = "  +0.84%  "
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "56.45"
$cell.Style = $origStyle

$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  +3.01%  "
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("E47").Value = "  +1.97%  "
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.413"
$cell.Style = $origStyle

$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("E50").Value = "  -1.53%  "
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.51"
$cell.Style = $origStyle

$ws.Range("E51").Value = "  -0.35%  "
